# Update model holdings "Weight" (D) and "Percent Change" (E) figures
# for rows 2-38, and refresh the "as of" date in the confidential
# disclosure footer (A41) from 2021-04-23 to 2021-04-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships protected; unprotect to allow writes, then
# re-protect once all edits are applied.
$ws.Unprotect()

$dValues = @{}   # column D ("Weight") new values, by row
$eValues = @{}   # column E ("Percent Change") new values, by row

$dValues[2] = 0.03027394334750616
$eValues[2] = 0.0000897827258035111
$dValues[3] = 0.03004484827305027
$eValues[3] = -0.0002714019851115523
$dValues[4] = 0.03052866600656222
$eValues[4] = 0.0320202996616723
$dValues[5] = 0.06486264002951056
$eValues[5] = 0.02038983740810796
$dValues[6] = 0.01538936455239094
$eValues[6] = -0.01203542502460064
$dValues[7] = 0.0157761081018114
$eValues[7] = 0.02606512589529641
$dValues[8] = 0.02994117304444056
$eValues[8] = -0.005874800606932906
$dValues[9] = 0.03386724134583649
$eValues[9] = -0.001490483833982981
$dValues[10] = 0.02944551226046945
$eValues[10] = 0.009395707645138973
$dValues[11] = 0.03154056730153183
$eValues[11] = 0.002671492588762314
$dValues[12] = 0.01348748713767061
$eValues[12] = 0.03303584280984584
$dValues[13] = 0.01480264818374202
$eValues[13] = -0.002754315093646764
$dValues[14] = 0.0163152581202555
$eValues[14] = 0.009079550187421903
$dValues[15] = 0.007986487273574466
$eValues[15] = 0.03257487359004263
$dValues[16] = 0.00710816005591473
$eValues[16] = 0.01941986234021642
$dValues[17] = 0.03173490981808298
$eValues[17] = 0.008350820093357969
$dValues[18] = 0.03000135903857728
$eValues[18] = 0.0008542141230067912
$dValues[19] = 0.0311637253188887
$eValues[19] = 0.02560508363704317
$dValues[20] = 0.02923194905546819
$eValues[20] = 0.006342775545445489
$dValues[21] = 0.04495117168540199
$eValues[21] = 0.00494104435710252
$dValues[22] = 0.03294212437144466
$eValues[22] = 0.0122881980256373
$dValues[23] = 0.03144135748539034
$eValues[23] = -0.01204112507332722
$dValues[24] = 0.02956491350690197
$eValues[24] = -0.008510638297872353
$dValues[25] = 0.01499388432640224
$eValues[25] = 0.01316862836499255
$dValues[26] = 0.0144621119459491
$eValues[26] = 0.006846556584776398
$dValues[27] = 0.03041780728832974
$eValues[27] = -0.006038053780804575
$dValues[28] = 0.03045799600054361
$eValues[28] = -0.0228454869964303
$dValues[29] = 0.03042110781058885
$eValues[29] = 0.001531686770055662
$dValues[30] = 0.02780379365911429
$eValues[30] = 0.01083032490974722
$dValues[31] = 0.03556468052886015
$eValues[31] = 0.01393688278934202
$dValues[32] = 0.03114392218533403
$eValues[32] = -0.01693753000068576
$dValues[33] = 0.02994389112159512
$eValues[33] = 0.02427511800404591
$dValues[34] = 0.03108781330692915
$eValues[34] = -0.01111638480177857
$dValues[35] = 0.03014502883103268
$eValues[35] = -0.0002318571759797328
$dValues[36] = 0.02915914341739958
$eValues[36] = 0.00246354617484501
$dValues[37] = 0.03199720426349816
$eValues[37] = -0.01155283724091061
$dValues[38] = 0.9999999999999999
$eValues[38] = 0.004686353311200264

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}

# Update the confidential disclosure text with the new "as of" date
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."
$ws.Range("A41").Value = $newText

# Restore sheet protection to match the original workbook state
$ws.Protect()
